# Updated iOS Tasks Status.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("iOS_Estimate")

# --- H27: turn the comment into rich text - plain run + a bold, appended sentence ---
# (written before H14 so the shared-string table keeps the same slot order as the original file)
$h27 = $ws.Range("H27")
$h27.Value = "Set up the hybris server and done initialisation to set up the code.Development Done.Not Getting response from hybris sever. Got the Successful response on 7th Aug."
$boldRun = $h27.Characters(125, 40)
$boldRun.Font.Bold = $true

# --- H14: append extra status detail to the existing comment text ---
$ws.Range("H14").Value = "Started Writing client code and methods.Waiting for updated webservice from Swarnima.Development Done.Getting 200 OK but not getting desired response i.e. ProfilepicURL."

# --- I15: new status ---
$ws.Range("I15").Value = "In progress"

# --- Row heights ---
$ws.Range("A14").RowHeight = 70
$ws.Range("A27").RowHeight = 56

# --- Row 22: start date + dependency comment + new "in progress" status ---
$ws.Range("F22").Value = 42226
$ws.Range("F22").NumberFormat = "d-mmm"
$ws.Range("H22").Value = "Need Webervice from Swarnima to proceed."
$ws.Range("I22").Value = "in progress"

# --- Row 23: start date + status ---
$ws.Range("F23").Value = 42226
$ws.Range("F23").NumberFormat = "d-mmm"
$ws.Range("I23").Value = "In progress"

# --- Row 26: end date + status -> completed ---
$ws.Range("G26").Value = 42223
$ws.Range("G26").NumberFormat = "d-mmm"
$ws.Range("I26").Value = "completed"

# --- Row 27: end date + status -> completed ---
$ws.Range("G27").Value = 42223
$ws.Range("G27").NumberFormat = "d-mmm"
$ws.Range("I27").Value = "completed"

# --- Selection moved to I16 ---
$ws.Range("I16").Select()
